$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.039628898182855
$ws.Range("D2").Value = 1.040077909014349
$ws.Range("E2").Value = 1.03806582671003
$ws.Range("F2").Value = 1.047367141158597
$ws.Range("I2").Value = 1.038048715023419
$ws.Range("J2").Value = 1.044719881389427
$ws.Range("K2").Value = 1.042861279582439
$ws.Range("L2").Value = 1.04085492065352
$ws.Range("M2").Value = 1.050129977221782
$ws.Range("N2").Value = 1.046203503452045

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.041138098692954
$ws.Range("D3").Value = 1.041421449347817
$ws.Range("E3").Value = 1.039368302121845
$ws.Range("F3").Value = 1.048991841193356
$ws.Range("I3").Value = 1.038444277576797
$ws.Range("J3").Value = 1.045871568733474
$ws.Range("K3").Value = 1.044013766300542
$ws.Range("L3").Value = 1.041966031120062
$ws.Range("M3").Value = 1.051564400405228
$ws.Range("N3").Value = 1.047356826324221

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.042112650351586
$ws.Range("D4").Value = 1.042289153832775
$ws.Range("E4").Value = 1.040209523305763
$ws.Range("F4").Value = 1.050041275052036
$ws.Range("I4").Value = 1.0386975721454
$ws.Range("J4").Value = 1.046614419928621
$ws.Range("K4").Value = 1.044757317260793
$ws.Range("L4").Value = 1.04268288895571
$ws.Range("M4").Value = 1.052490241827758
$ws.Range("N4").Value = 1.048100732453272

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.042521881601634
$ws.Range("D5").Value = 1.042653547613996
$ws.Range("E5").Value = 1.040562804346008
$ws.Range("F5").Value = 1.05048202253525
$ws.Range("I5").Value = 1.038803422983074
$ws.Range("J5").Value = 1.046926154563227
$ws.Range("K5").Value = 1.045069389663082
$ws.Range("L5").Value = 1.042983758457948
$ws.Range("M5").Value = 1.052878917450248
$ws.Range("N5").Value = 1.048412909786815

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.042590566013208
$ws.Range("D6").Value = 1.042714708321079
$ws.Range("E6").Value = 1.04062210038727
$ws.Range("F6").Value = 1.050556000779053
$ws.Range("I6").Value = 1.03882115869598
$ws.Range("J6").Value = 1.046978463537032
$ws.Range("K6").Value = 1.045121757902395
$ws.Range("L6").Value = 1.043034246780628
$ws.Range("M6").Value = 1.052944145910667
$ws.Range("N6").Value = 1.048465293045363

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.0421181203563
$ws.Range("D7").Value = 1.042294024400607
$ws.Range("E7").Value = 1.040214245303137
$ws.Range("F7").Value = 1.050047166036801
$ws.Range("I7").Value = 1.038698989018004
$ws.Range("J7").Value = 1.046618587534176
$ws.Range("K7").Value = 1.044761489208479
$ws.Range("L7").Value = 1.042686911135819
$ws.Range("M7").Value = 1.052495437473317
$ws.Range("N7").Value = 1.048104905977304

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.040139359529953
$ws.Range("D8").Value = 1.040532312017948
$ws.Range("E8").Value = 1.03850633311389
$ws.Range("F8").Value = 1.047916606266406
$ws.Range("I8").Value = 1.038182949890783
$ws.Range("J8").Value = 1.04510959383018
$ws.Range("K8").Value = 1.043251224318155
$ws.Range("L8").Value = 1.04123086534082
$ws.Range("M8").Value = 1.050615234248282
$ws.Range("N8").Value = 1.046593769329164

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.03663679157086
$ws.Range("D9").Value = 1.037414940365798
$ws.Range("E9").Value = 1.035484456848852
$ws.Range("F9").Value = 1.044147636363842
$ws.Range("I9").Value = 1.037253116717917
$ws.Range("J9").Value = 1.042432116002652
$ws.Range("K9").Value = 1.040572909324265
$ws.Range("L9").Value = 1.038648717312256
$ws.Range("M9").Value = 1.04728386258752
$ws.Range("N9").Value = 1.043912489176011

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.034290594350215
$ws.Range("D10").Value = 1.03532748602872
$ws.Range("E10").Value = 1.033461153057041
$ws.Range("F10").Value = 1.041624522195846
$ws.Range("I10").Value = 1.036619260395995
$ws.Range("J10").Value = 1.040634303640133
$ws.Range("K10").Value = 1.038775489138307
$ws.Range("L10").Value = 1.036915847224548
$ws.Range("M10").Value = 1.045050144650138
$ws.Range("N10").Value = 1.042112123713725

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.033271892256465
$ws.Range("D11").Value = 1.034421308401347
$ws.Range("E11").Value = 1.032582873595427
$ws.Range("F11").Value = 1.040529364994142
$ws.Range("I11").Value = 1.036341439915583
$ws.Range("J11").Value = 1.0398526919738
$ws.Range("K11").Value = 1.037994274312641
$ws.Range("L11").Value = 1.036162691405891
$ws.Range("M11").Value = 1.044079758257158
$ws.Range("N11").Value = 1.041329402069173

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.032893071698054
$ws.Range("D12").Value = 1.034084359731856
$ws.Range("E12").Value = 1.032256306073205
$ws.Range("F12").Value = 1.040122167404653
$ws.Range("I12").Value = 1.036237737140589
$ws.Range("J12").Value = 1.039561885363288
$ws.Range("K12").Value = 1.037703649360021
$ws.Range("L12").Value = 1.035882505452889
$ws.Range("M12").Value = 1.043718826143631
$ws.Range("N12").Value = 1.041038182479919

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.032974349597185
$ws.Range("D13").Value = 1.034156652527461
$ws.Range("E13").Value = 1.032326371232429
$ws.Range("F13").Value = 1.04020953130135
$ws.Range("I13").Value = 1.036260004778682
$ws.Range("J13").Value = 1.039624286281557
$ws.Range("K13").Value = 1.037766009759642
$ws.Range("L13").Value = 1.035942625893204
$ws.Range("M13").Value = 1.043796269538838
$ws.Range("N13").Value = 1.04110067201465

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.033240587651923
$ws.Range("D14").Value = 1.034393463391992
$ws.Range("E14").Value = 1.032555886314599
$ws.Range("F14").Value = 1.040495714295891
$ws.Range("I14").Value = 1.036332878193401
$ws.Range("J14").Value = 1.039828663674782
$ws.Range("K14").Value = 1.037970260329906
$ws.Range("L14").Value = 1.036139539972724
$ws.Range("M14").Value = 1.044049933502178
$ws.Range("N14").Value = 1.041305339647216

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.033404568491802
$ws.Range("D15").Value = 1.034539323209162
$ws.Range("E15").Value = 1.032697253426853
$ws.Range("F15").Value = 1.040671986733054
$ws.Range("I15").Value = 1.036377710502054
$ws.Range("J15").Value = 1.03995452331002
$ws.Range("K15").Value = 1.038096046368669
$ws.Range("L15").Value = 1.036260807972289
$ws.Range("M15").Value = 1.044206159329468
$ws.Range("N15").Value = 1.041431378017571

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.034358142669236
$ws.Range("D16").Value = 1.035387576911821
$ws.Range("E16").Value = 1.033519395005988
$ws.Range("F16").Value = 1.041697147695425
$ws.Range("I16").Value = 1.036637627402078
$ws.Range("J16").Value = 1.040686109610066
$ws.Range("K16").Value = 1.038827273566761
$ws.Range("L16").Value = 1.036965771782596
$ws.Range("M16").Value = 1.045114478247366
$ws.Range("N16").Value = 1.042164003254078

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.034955541610273
$ws.Range("D17").Value = 1.035919042995159
$ws.Range("E17").Value = 1.034034514211406
$ws.Range("F17").Value = 1.042339490956214
$ws.Range("I17").Value = 1.036799765244197
$ws.Range("J17").Value = 1.041144165889167
$ws.Range("K17").Value = 1.039285165473765
$ws.Range("L17").Value = 1.037407218625433
$ws.Range("M17").Value = 1.045683385884701
$ws.Range("N17").Value = 1.042622710025634

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.035303726520786
$ws.Range("D18").Value = 1.036228817238996
$ws.Range("E18").Value = 1.034334765201651
$ws.Range("F18").Value = 1.042713905768171
$ws.Range("I18").Value = 1.036894013895541
$ws.Range("J18").Value = 1.041411039296425
$ws.Range("K18").Value = 1.039551964900247
$ws.Range("L18").Value = 1.037664436358197
$ws.Range("M18").Value = 1.046014914793531
$ws.Range("N18").Value = 1.042889962423742

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.035422403490399
$ws.Range("D19").Value = 1.036334405047426
$ws.Range("E19").Value = 1.034437107782444
$ws.Range("F19").Value = 1.042841528966636
$ws.Range("I19").Value = 1.036926095468538
$ws.Range("J19").Value = 1.041501985079137
$ws.Range("K19").Value = 1.039642889166788
$ws.Range("L19").Value = 1.03775209527403
$ws.Range("M19").Value = 1.046127906128574
$ws.Range("N19").Value = 1.042981037359897

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.034891474136126
$ws.Range("D20").Value = 1.035862044631104
$ws.Range("E20").Value = 1.033979268509428
$ws.Range("F20").Value = 1.042270599874106
$ws.Range("I20").Value = 1.036782402886625
$ws.Range("J20").Value = 1.041095052194421
$ws.Range("K20").Value = 1.039236067148873
$ws.Range("L20").Value = 1.037359883618521
$ws.Range("M20").Value = 1.04562237913564
$ws.Range("N20").Value = 1.042573526583809

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.033162199119451
$ws.Range("D21").Value = 1.034323738338107
$ws.Range("E21").Value = 1.032488309172012
$ws.Range("F21").Value = 1.040411451881257
$ws.Range("I21").Value = 1.036311432838148
$ws.Range("J21").Value = 1.039768492979023
$ws.Range("K21").Value = 1.037910126033597
$ws.Range("L21").Value = 1.036081565628034
$ws.Range("M21").Value = 1.043975249329894
$ws.Range("N21").Value = 1.041245083502166

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.032072446242206
$ws.Range("D22").Value = 1.03335449175172
$ws.Range("E22").Value = 1.031548938891877
$ws.Range("F22").Value = 1.039240168384233
$ws.Range("I22").Value = 1.036012375033143
$ws.Range("J22").Value = 1.038931644029533
$ws.Range("K22").Value = 1.037073863963124
$ws.Range("L22").Value = 1.035275342528638
$ws.Range("M22").Value = 1.042936808032848
$ws.Range("N22").Value = 1.04040704613118

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.03265038452436
$ws.Range("D23").Value = 1.033868505267696
$ws.Range("E23").Value = 1.032047104093306
$ws.Range("F23").Value = 1.039861315981532
$ws.Range("I23").Value = 1.03617119114768
$ws.Range("J23").Value = 1.039375540645457
$ws.Range("K23").Value = 1.037517430621764
$ws.Range("L23").Value = 1.035702975702385
$ws.Range("M23").Value = 1.043487576886709
$ws.Range("N23").Value = 1.040851573131215

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.034920424297489
$ws.Range("D24").Value = 1.035887800421178
$ws.Range("E24").Value = 1.034004232311101
$ws.Range("F24").Value = 1.042301729571031
$ws.Range("I24").Value = 1.036790249188658
$ws.Range("J24").Value = 1.041117245495596
$ws.Range("K24").Value = 1.039258253437423
$ws.Range("L24").Value = 1.037381273107336
$ws.Range("M24").Value = 1.045649946402963
$ws.Range("N24").Value = 1.042595751402017

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.037544212837484
$ws.Range("D25").Value = 1.038222444414285
$ws.Range("E25").Value = 1.036267188797742
$ws.Range("F25").Value = 1.045123804846854
$ws.Range("I25").Value = 1.037495949408217
$ws.Range("J25").Value = 1.043126537269007
$ws.Range("K25").Value = 1.041267380054337
$ws.Range("L25").Value = 1.039318250563242
$ws.Range("M25").Value = 1.048147317504886
$ws.Range("N25").Value = 1.044607896600215

